$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp footer text (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 23 de Mayo de 2020 a las 00:05"

# Refresh country statistics (new COVID-19 snapshot values)

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1642717
$ws.Range("C4").Value = 21815
$ws.Range("D4").Value = 396139
$ws.Range("E4").Value = 1149039
$ws.Range("G4").Value = 1185
$ws.Range("H4").Value = 97539

# Colombia (row 39)
$ws.Range("B39").Value = 19131
$ws.Range("C39").Value = 801
$ws.Range("D39").Value = 4575
$ws.Range("E39").Value = 13874
$ws.Range("G39").Value = 30
$ws.Range("H39").Value = 682

# El Salvador (row 92)
$ws.Range("B92").Value = 1728
$ws.Range("C92").Value = 161
$ws.Range("D92").Value = 402
$ws.Range("E92").Value = 1314
$ws.Range("H92").Value = 12

# Lituania (row 93)
$ws.Range("B93").Value = 1725
$ws.Range("C93").Value = 85
$ws.Range("D93").Value = 562
$ws.Range("E93").Value = 1130
$ws.Range("H93").Value = 33

# Somalia (row 94)
$ws.Range("B94").Value = 1604
$ws.Range("C94").Value = 11
$ws.Range("D94").Value = 1111
$ws.Range("E94").Value = 432

# Gabon (row 95)
$ws.Range("B95").Value = 1594
$ws.Range("D95").Value = 204
$ws.Range("E95").Value = 1329
$ws.Range("H95").Value = 61

# Guinea-Bisau (row 103)
$ws.Range("B103").Value = 1114
$ws.Range("C103").Value = 5
$ws.Range("E103").Value = 1066

# Mozambique (row 157)
$ws.Range("B157").Value = 175
$ws.Range("C157").Value = 15
$ws.Range("D157").Value = 68
$ws.Range("E157").Value = 107

# Uganda (row 158)
$ws.Range("B158").Value = 164
$ws.Range("C158").Value = 2
$ws.Range("D158").Value = 48
$ws.Range("E158").Value = 116

# Guyana (row 165)
$ws.Range("B165").Value = 128
$ws.Range("C165").Value = 3
$ws.Range("D165").Value = 81
$ws.Range("E165").Value = 38
$ws.Range("H165").Value = 9

# Bermudas (row 166)
$ws.Range("B166").Value = 127
$ws.Range("D166").Value = 57
$ws.Range("E166").Value = 60
$ws.Range("H166").Value = 10
